# Reshape of comparison tables
# ------------------------------------------------------------------
# The "BSS" rows of both comparison blocks (FRESH_LS / DRY_LS) were
# previously blank placeholders; the author filled them in with the
# actual BSS estimates, which ripples into the "Delta" (ABS) rows
# below them. At the same time the numeric columns were reformatted:
# the variance columns (C, D) now share a single 4-decimal format and
# the two small-magnitude columns (E, F) were switched to scientific
# notation.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty BSS values for both comparison blocks.
# (The "Delta" rows below each, e.g. C4/D4 and C7/D7, hold
# =ABS(...) formulas referencing these cells, so they recalc
# automatically once the values land.)
$ws.Range("C3").Value = 0.18285151213991993
$ws.Range("D3").Value = 0.26441975235277537
$ws.Range("C6").Value = 4.5602399335886421
$ws.Range("D6").Value = 4.9244855632738407

# Reformat the numeric columns of the reshaped tables.
$ws.Range("C2:D7").NumberFormat = "0.0000"
$ws.Range("E2:F7").NumberFormat = "0.000E+00"

# Reflect the reshaped table in the sheet's current selection.
[void]$ws.Range("A1:F7").Select()
